$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1 - Product Burndown")

# "Fächer sind zu den Tutoren gespeichert" -> record the missing "Actual Hours"
# entry for the sprint starting 41419 and correct the planned-hours entry for
# the following sprint (41426).
$ws.Range("E7").Value = 10
$ws.Range("B8").Value = 7

# Leave the selection where the author last left it before saving.
$ws.Activate()
[void]$ws.Range("N19").Select()
